$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Header row: add two new columns (mean/std) for each of the three groups
# ---------------------------------------------------------------------------
$ws.Range("B1").Value = "Algorithm"
$ws.Range("C1").Value = "State Based mean"
$ws.Range("D1").Value = "State Based std"
$ws.Range("E1").Value = "Non State mean"
$ws.Range("F1").Value = "Non State std"
$ws.Range("G1").Value = "One Sided mean"
$ws.Range("H1").Value = "One Sided std"

# Give the two new header cells (G1, H1) the same look as the existing
# header cells (bold, centered, bordered) by copying the format from D1.
$ws.Range("D1").Copy() | Out-Null
$ws.Range("F1:H1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# ---------------------------------------------------------------------------
# 2. Body rows: 7 algorithms (CART renamed to DTREE, NB removed), each with
#    six numeric columns (mean/std per group) instead of three.
# ---------------------------------------------------------------------------
$names = @("LR", "LDA", "KNN", "DTREE", "RTREE", "XTREE", "SVM")

$data = @(
    @(0.7605902407559192, 0.02912092038057533, 0.5769184139509725, 0.03116412515734985, 0.7425940608799981, 0.02060120262929459),
    @(0.7576326266171074, 0.02758041201676718, 0.5476353479174623, 0.03752369141805782, 0.7267392633558764, 0.01911661650980238),
    @(0.7863669187482778, 0.02754914785742697, 0.6713884126468413, 0.03660659015393389, 0.7466812456377306, 0.01730827436236156),
    @(0.7554103753930235, 0.02635369798642654, 0.6665567863942251, 0.03908385819458696, 0.7220109420252856, 0.01941248095376418),
    @(0.7599349068026771, 0.02606021869132721, 0.5416935783742984, 0.03385873008535697, 0.7321594678086751, 0.01966134238044209),
    @(0.8208539601169849, 0.02808521397824066, 0.7067501603665434, 0.03524409152271004, 0.7880615701290051, 0.02314158065942502),
    @(0.800511415845199,  0.02721022775054627, 0.6718074731601458, 0.03072141214940161, 0.7713099607472109, 0.01254476845599847)
)

for ($i = 0; $i -lt $names.Length; $i++) {
    $row = $i + 2

    $ws.Cells.Item($row, 1).Value = $i
    $ws.Cells.Item($row, 2).Value = $names[$i]

    $values = $data[$i]
    for ($j = 0; $j -lt $values.Length; $j++) {
        $ws.Cells.Item($row, 3 + $j).Value = $values[$j]
    }
}

# ---------------------------------------------------------------------------
# 3. Remove the now-unused last row (previously row 9, the SVM entry) since
#    the table now only has 7 data rows (rows 2-8).
# ---------------------------------------------------------------------------
$ws.Rows.Item(9).Delete() | Out-Null

$ws.Range("A1").Select() | Out-Null
